$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column G ("K" - strikeout count) values for rows 2-6, as part of
# regenerating save_data to use K instead of Strike# (recomputed stat).
$ws.Range("G2").Value = 9
$ws.Range("G3").Value = 11
$ws.Range("G4").Value = 7
$ws.Range("G5").Value = 8
$ws.Range("G6").Value = 4
